$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 2 with the new agenda entry
$ws.Range("A2").Value = "Giovani"
$ws.Range("B2").Value = "0701"
$ws.Range("C2").Value = "Usina Amaral"
$ws.Range("D2").Value = "Local sem comunicação de alarmes e câmeras, disse o Alisson que a energia lá foi normalizada."
$ws.Range("E2").Value = ""
$ws.Range("F2").Value = ""
$ws.Range("G2").Value = "Pendente"
$ws.Range("H2").Value = "Maxvel: 31 / Forte: 11"
$ws.Range("I2").Value = ""

# Clear out rows 3 through 10 (values only, keep styles); row height reverts to default automatically
$ws.Range("A3:I10").ClearContents()

# Change the top-left visible cell of the sheet view from D1 to E1
$ws.Application.ActiveWindow.ScrollColumn = 5

Write-Host "Done"
